$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the extra "{{ clients[0].name}} " paragraph that duplicated the
#    greeting placeholder right after the date line.
# ---------------------------------------------------------------------------
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Text -eq "{{ clients[0].name}} `r") {
        $pp.Range.Delete()
        $found = $true
        break
    }
}
Write-Host "Removed duplicate clients[0].name paragraph:" $found

# ---------------------------------------------------------------------------
# 2. "Sincerely yours," moves down one paragraph: the paragraph that used to
#    hold it becomes blank, and the paragraph that used to hold
#    "{{ attorneys[0].signature }}" now reads "Sincerely yours,".
# ---------------------------------------------------------------------------
$ok1 = $d.Content.Find.Execute("Sincerely yours,", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
Write-Host "Cleared old Sincerely yours,:" $ok1

$ok2 = $d.Content.Find.Execute("{{ attorneys[0].signature }}", $true, $false, $false, $false, $false, $true, 1, $false, "Sincerely yours,", 2)
Write-Host "Inserted Sincerely yours, on signature line:" $ok2

# ---------------------------------------------------------------------------
# 3. "Staff Attorney" is replaced by the attorney signature merge-field.
# ---------------------------------------------------------------------------
$ok3 = $d.Content.Find.Execute("Staff Attorney", $true, $false, $false, $false, $false, $true, 1, $false, "{{ attorneys[0].signature }}", 2)
Write-Host "Replaced Staff Attorney with signature field:" $ok3

# ---------------------------------------------------------------------------
# 4. Insert a new bold "Advocate Title" placeholder line right after the
#    signature placeholder paragraph.
# ---------------------------------------------------------------------------
$sigPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Text -eq "{{ attorneys[0].signature }}`r") {
        $sigPara = $pp
    }
}
$sigPara.Range.InsertParagraphAfter()

$newIndex = $sigPara.Index + 1
$titlePara = $d.Paragraphs.Item($newIndex)
$tr = $titlePara.Range
$tr2 = $d.Range($tr.Start, $tr.End - 1)
$dash = [char]0x2013
$tr2.Text = "Advocate Title " + $dash + " need variable"
$tr2.Font.Bold = 1
Write-Host "New title paragraph:" $d.Paragraphs.Item($newIndex).Range.Text

# ---------------------------------------------------------------------------
# 5. "Tel.:" becomes "Direct Telephone:" in front of the phone placeholder.
# ---------------------------------------------------------------------------
$ok4 = $d.Content.Find.Execute("Tel.:", $true, $false, $false, $false, $false, $true, 1, $false, "Direct Telephone:", 2)
Write-Host "Relabeled phone line:" $ok4

# ---------------------------------------------------------------------------
# 6. Merge the two footer runs that together spell out the office address so
#    they become a single run (no textual change, just a run-merge).
# ---------------------------------------------------------------------------
$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(2)
$ok5 = $ftr.Range.Find.Execute("Boston, MA 02114", $true, $false, $false, $false, $false, $true, 1, $false, "Boston, MA 02114", 2)
Write-Host "Merged footer address runs:" $ok5
